$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.611.01"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.640.25"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "214.86"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "0.0625"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "19.06"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.868.59"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "1.638.08"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "64.79"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "26.621.74"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "215.53"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "6.26"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +12.77%  "
$ws.Range("D25").Value = "145.12"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "7.11"
$ws.Range("E28").Value = "  +3.93%  "
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D34").Value = "1.274.30"
$ws.Range("E34").Value = "  +4.56%  "
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "0.0178"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("D38").Value = "0.533"
$ws.Range("E38").Value = "  +6.37%  "
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "1.779.48"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").Value = "91.31"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.00"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "0.0515"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("D49").Value = "7.72"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  +0.13%  "
